$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Rename a couple of helper-column CQL ids (C11 / C12) - content only
#    renamed, same position, same style.
# ---------------------------------------------------------------------------
$ws.Range("C11").Value = "H-HTNStage2LastBPSetOffice"
$ws.Range("C12").Value = "H-HTNStage2AverageBPOffice"

# ---------------------------------------------------------------------------
# 2) Re-lay rows 13..17 out as rows 14..18 (shift down by one row), then
#    turn old row 13 ("Consider HTN Stage 2") into two rows: the header
#    row 13 (unchanged A/B, new C) and a brand-new detail row 14.
#
#    We copy formats from the row that currently holds the right look,
#    using column runs that exactly match the destination's non-blank
#    cells so we don't leave stray styled-blank cells behind (mirrors
#    how the destination workbook actually lays these rows out).
#    Processed bottom-up so each format source is read before it is
#    itself overwritten.
# ---------------------------------------------------------------------------

# row 18  <=  old row 17 ("Prescribe Ambulatory BP Monitoring")
$ws.Range("A17:F17").Copy()
$ws.Range("A18:F18").PasteSpecial(-4122)
$ws.Range("H17:J17").Copy()
$ws.Range("H18:J18").PasteSpecial(-4122)
$ws.Range("M17:O17").Copy()
$ws.Range("M18:O18").PasteSpecial(-4122)

$ws.Range("A18").Value = "No Further Action"
$ws.Range("B18").Value = "None"
$ws.Range("C18").Value = "H-NoFurtherAction"
$ws.Range("D18").Value = "18-80"
$ws.Range("E18").Value = "130-179"
$ws.Range("F18").Value = "80-119"
$ws.Range("H18").Value = $true
$ws.Range("I18").Value = $false
$ws.Range("J18").Value = $false
$ws.Range("M18").Value = $false
$ws.Range("N18").Value = $false
$ws.Range("O18").Value = $false

# row 17  <=  old row 16 ("Prescribe HBP or ABP Monitoring")
$ws.Range("A16:F16").Copy()
$ws.Range("A17:F17").PasteSpecial(-4122)
$ws.Range("H16:J16").Copy()
$ws.Range("H17:J17").PasteSpecial(-4122)
$ws.Range("M16:N16").Copy()
$ws.Range("M17:N17").PasteSpecial(-4122)

$ws.Range("A17").Value = "Prescribe Ambulatory BP Monitoring"
$ws.Range("B17").Value = "Prescribe Ambulatory BP Monitoring"
$ws.Range("C17").Value = "H-PrescribeAmbulatoryBPMonitoring"
$ws.Range("D17").Value = "18-80"
$ws.Range("E17").Value = "130-179"
$ws.Range("F17").Value = "80-119"
$ws.Range("H17").Value = $true
$ws.Range("I17").Value = $false
$ws.Range("J17").Value = $false
$ws.Range("M17").Value = $false
$ws.Range("N17").Value = $false
$ws.Range("O17").Value = $true

# row 16  <=  old row 15 ("Consider HTN Stage 1")
$ws.Range("A15:F15").Copy()
$ws.Range("A16:F16").PasteSpecial(-4122)
$ws.Range("H15:J15").Copy()
$ws.Range("H16:J16").PasteSpecial(-4122)
$ws.Range("M15:M15").Copy()
$ws.Range("M16:M16").PasteSpecial(-4122)

$ws.Range("A16").Value = "Prescribe HBP or ABP Monitoring"
$ws.Range("B16").Value = "Prescribe HBP or ABP Monitoring"
$ws.Range("C16").Value = "H-PrescribeHBPABPMonitoring"
$ws.Range("D16").Value = "18-80"
$ws.Range("E16").Value = "130-179"
$ws.Range("F16").Value = "80-119"
$ws.Range("H16").Value = $true
$ws.Range("I16").Value = $false
$ws.Range("J16").Value = $false
$ws.Range("M16").Value = $false
$ws.Range("N16").Value = $true

# row 15  <=  old row 14 ("Consider HTN Stage 1" source formats)
$ws.Range("A14:F14").Copy()
$ws.Range("A15:F15").PasteSpecial(-4122)
$ws.Range("H14:J14").Copy()
$ws.Range("H15:J15").PasteSpecial(-4122)

$ws.Range("A15").Value = "Consider HTN Stage 1"
$ws.Range("B15").Value = "Consider HTN Stage 1"
$ws.Range("C15").Value = "H-ConsiderHTNStage1"
$ws.Range("D15").Value = "18-80"
$ws.Range("E15").Value = "130-179"
$ws.Range("F15").Value = "80-119"
$ws.Range("H15").Value = $true
$ws.Range("I15").Value = $false
$ws.Range("J15").Value = $false
$ws.Range("M15").Value = $true

# row 14 (brand-new detail row) <= old row 13's C:L formats (style 20)
$ws.Range("C13:F13").Copy()
$ws.Range("C14:F14").PasteSpecial(-4122)
$ws.Range("H13:L13").Copy()
$ws.Range("H14:L14").PasteSpecial(-4122)

$ws.Range("C14").Value = "H-ConsiderHTNStage2"
$ws.Range("D14").Value = "18-80"
$ws.Range("E14").Value = "130-179"
$ws.Range("F14").Value = "80-119"
$ws.Range("H14").Value = $true
$ws.Range("I14").Value = $true
$ws.Range("J14").Value = $false
$ws.Range("K14").Value = $false
$ws.Range("L14").Value = $false

# row 13 stays "Consider HTN Stage 2" (A/B untouched) but its C cell now
# gets the divider look (style 7, sourced from the old C12 cell) and a new id
$ws.Range("C12").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value = "H-ConsiderHTNStage2LastBPSetOffice"

# ---------------------------------------------------------------------------
# 3) Selection cursor moved (matches the captured UI state in the diff)
# ---------------------------------------------------------------------------
$ws.Range("C21").Select()
